$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 364
$ws1.Range("F4").Value = 2980
$ws1.Range("F5").Value = 75
$ws1.Range("F6").Value = 622

# Sheet "全部类型" (All types) - update the same values for the corresponding rows
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F5").Value = 364
$ws2.Range("F6").Value = 2980
$ws2.Range("F7").Value = 75
$ws2.Range("F8").Value = 622
